$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 90

$ws.Cells.Item($row, 1).Value = 90
$ws.Cells.Item($row, 2).Value = "Conhecimentos Específicos"
$ws.Cells.Item($row, 3).Value = "Gestão de Estoques"
$ws.Cells.Item($row, 4).Value = "<b>VMI</b>:`n<i>Conceito</i>"
$ws.Cells.Item($row, 5).Value = "modelo de reposição de estoques no qual o fornecedor assume a responsabilidade de monitorar e reabastecer o estoque do cliente com base na demanda real e nos níveis de estoque disponíveis."
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0

$ws.Rows.Item($row).AutoFit()
